# Append the newest profit-allocation row (run date 2025-10-10) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to treat the value as literal text instead
# of auto-converting the "MM/DD/YYYY" looking string into a date serial,
# matching how the other Date column cells are stored (t="inlineStr"/"s").
$ws.Range("A39").Value = "'10/10/2025"
# Reset to the default "Normal" style so no date/quote-prefix number format
# (and no extra style record) gets attached to the new cell, consistent
# with the unstyled cells used by every other data row.
$ws.Range("A39").Style = "Normal"

$ws.Range("B39").Value = 0.1531112533627095
$ws.Range("C39").Value = 0.8468887466372905
